# Apply the "Added periodic & upfront related scenarios" edit:
#  - Change the repayment strategy value on ProductLoanInput!B17 from
#    "Mifos style" to "Penalties, Fees, Interest, Principal order" and
#    left/top-align that cell.
#  - Make ProductLoanInput the active sheet/tab (selecting B17), and
#    ProductLoanOutput no longer the active/selected sheet.

$wb   = $excel.ActiveWorkbook
$wsIn = $wb.Worksheets.Item("ProductLoanInput")

# Update the repayment strategy cell value and alignment.
$rng = $wsIn.Range("B17")
$rng.Value = "Penalties, Fees, Interest, Principal order"
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4160

# Make ProductLoanInput the active sheet and select B17 on it, which
# leaves ProductLoanOutput's selection untouched (still B1) but no
# longer the active/tabSelected sheet.
$wsIn.Activate() | Out-Null
$wsIn.Range("B17").Select() | Out-Null
